$d = $word.ActiveDocument

function Replace-TextPreservingFormat($doc, [string]$old, [string]$new) {
    # Plain Find.Execute (or Range.Text=) that matches/replaces the *entire*
    # text of a run resets that run's rPr to a generic default, losing
    # formatting (color, sz, u, etc.) that isn't shared with sibling runs.
    # Splitting the edit so the Find/Range op always leaves at least one
    # original character in place (a one-character "anchor" at the start)
    # keeps it a genuine partial-run edit and preserves the run's rPr.
    $rng = $doc.Content
    $found = $rng.Find.Execute($old)
    if (-not $found) {
        return $false
    }
    $startPos = $rng.Start
    $endPos = $rng.End

    if ($new.Length -eq 0) {
        # Full deletion: nothing to anchor on the replacement side, just
        # remove everything after the first character, then remove the
        # anchor character itself (now a true partial/final cleanup op).
        $sub = $doc.Range($startPos + 1, $endPos)
        $sub.Text = ""
        $lead = $doc.Range($startPos, $startPos + 1)
        $lead.Text = ""
        return $true
    }

    # Step 1: replace everything except the first character of the match
    # with everything except the first character of the replacement.
    $sub = $doc.Range($startPos + 1, $endPos)
    $sub.Text = $new.Substring(1)

    # Step 2: fix up the still-original leading anchor character (a
    # genuine partial, 1-of-N-character edit, so formatting is untouched).
    $lead = $doc.Range($startPos, $startPos + 1)
    $lead.Text = $new.Substring(0, 1)

    return $true
}

$changes = 0

if (Replace-TextPreservingFormat $d "Aplicativo de Controle Financeiro" "Gerenciador de Tarefas") { $changes++ }

if (Replace-TextPreservingFormat $d `
    "Implementação de um sistema de controle financeiro com Python e Firebase, incluindo integração de dados, autenticação, escalabilidade e testes automatizados." `
    "Plataforma desenvolvida para automação de integração de dados utilizando Python e serviços em nuvem, focada em escalabilidade e eficiência, com práticas de testes automatizados e segurança de dados.") { $changes++ }

if (Replace-TextPreservingFormat $d "https://github.com/fakeuser/finance-app" "https://github.com/fakeuser/task-manager") { $changes++ }

if (Replace-TextPreservingFormat $d `
    "Experiência em Python, manipulação de dados em bancos SQL e NoSQL, habilidades em cloud computing com AWS/Azure e integração de APIs. Competências relevantes: Python, AWS, Azure, SQL, Git, Machine Learning, APIs RESTful" `
    "Experiência sólida em desenvolvimento e integração de workflows escaláveis, deployment em cloud (AWS, Azure), codificação eficiente em Python. Competências: Python, AWS, Azure, SQL, Git") { $changes++ }

if (Replace-TextPreservingFormat $d `
    "Atuei como desenvolvedor full stack utilizando Python em aplicações escaláveis e modelagem de dados, com interface entre APIs e bancos de dados, além de projetos colaborativos focados em soluções eficientes." `
    "Atuação como Desenvolvedor Full Stack aplicando Python e integração de soluções em nuvem (AWS, Azure), com experiência em manipulação de dados, automação de sistemas e colaboração em equipes utilizando metodologias ágeis.") { $changes++ }

if (Replace-TextPreservingFormat $d `
    "Durante o bacharelado em Ciência da Computação, desenvolvi sólida base em programação Python, fundamentos de Machine Learning e análise de dados, além de experiência com projetos de mineração de dados e algoritmos avançados." `
    "Graduação em Ciência da Computação pela UFRJ, com sólida base em programação, estatística e projetos práticos voltados para o desenvolvimento em inteligência artificial e ciência de dados.") { $changes++ }

if (Replace-TextPreservingFormat $d `
    "O curso técnico em Mecatrônica aprimorou minha capacidade analítica e resolução de problemas utilizando tecnologias inovadoras e integração de sistemas automatizados." `
    "Curso Técnico em Engenharia Mecatrônica no SENAI com foco em automação, integração de sistemas e solução de problemas multidisciplinares.") { $changes++ }

Write-Output "Applied $changes/7 replacements"

$d.Save()
